# Apply cryptocurrency price/volume updates for Sun Oct 15 12:15:48 UTC 2023 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.044.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.565.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  +0.69%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.490"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.47%  "

$ws.Range("E7").Value = "  +0.65%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.248"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0595"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0861"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.789.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.570.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.041.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0703"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.81%  "

$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("E22").Value = "  +1.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.55%  "

$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("E28").Value = "  +1.22%  "

$ws.Range("E29").Value = "  +0.59%  "

$ws.Range("E30").Value = "  +4.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0472"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.425.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +18.09%  "

$ws.Range("E36").Value = "  +1.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.22%  "

$ws.Range("E38").Value = "  +1.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.531"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.63%  "

$ws.Range("E40").Value = "  +2.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.58%  "

$ws.Range("E43").Value = "  +0.50%  "

$ws.Range("E44").Value = "  +0.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.707.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.65%  "

$ws.Range("E49").Value = "  +2.09%  "

$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0961"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.45%  "
